$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B2").Value = "Beispielfirma GmbH"
$ws.Range("B3").Value = "Bf GmbH"
